# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1): rows 2-6 and 8 get bumped F values.
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 278
$ws1.Range("F3").Value = 172
$ws1.Range("F4").Value = 2062
$ws1.Range("F5").Value = 1650
$ws1.Range("F6").Value = 298
$ws1.Range("F8").Value = 673

# Sheet "全部类型" (index 4): same underlying rows (offset by the extra
# "演出" row), rows 2-6 and 9 get bumped F values.
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 278
$ws4.Range("F3").Value = 172
$ws4.Range("F4").Value = 2062
$ws4.Range("F5").Value = 1650
$ws4.Range("F6").Value = 298
$ws4.Range("F9").Value = 673
